# "Update countries & provincias Spain"
#
# The source "Pais" sheet is a COVID-19 snapshot: row 1 holds a last-updated
# timestamp, row 3 holds headers, and rows 4+ hold one country per row with
# columns B..H = Casos totales / Nuevos casos / Casos activos / Recuperados /
# Casos criticos / Muertes hoy / Muertes.
#
# This edit refreshes the snapshot timestamp, updates the numbers for several
# countries with newer totals, and fixes the row/label alignment for a few
# country pairs (Hong Kong/Tunez, Groenlandia/Islas Malvinas,
# Santa Sede/Islas Turcas y Caicos, Papua Nueva Guinea/Islas Virgenes
# Britanicas) whose names had drifted onto the wrong row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Snapshot timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Junio de 2020 a las 14:35"

# --- Updated case counts for existing countries ---
# Row 4: Estados Unidos
$ws.Range("B4").Value = 2183126
$ws.Range("C4").Value = 176
$ws.Range("E4").Value = 1174790
$ws.Range("G4").Value = 38
$ws.Range("H4").Value = 118321

# Row 26: Bielorrusia
$ws.Range("B26").Value = 55369
$ws.Range("C26").Value = 689
$ws.Range("D26").Value = 31273
$ws.Range("E26").Value = 23778
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = 318

# Row 29: Paises Bajos
$ws.Range("B29").Value = 49087
$ws.Range("C29").Value = 139
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = 6070

# Row 79: Consejo Danes para los Refugiados
$ws.Range("B79").Value = 4974
$ws.Range("C79").Value = 137
$ws.Range("D79").Value = 628
$ws.Range("E79").Value = 4234

# Row 89: Etiopia
$ws.Range("B89").Value = 3630
$ws.Range("C89").Value = 109
$ws.Range("D89").Value = 738
$ws.Range("E89").Value = 2831
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 61

# Row 100: Croacia
$ws.Range("B100").Value = 2255
$ws.Range("C100").Value = 1
$ws.Range("E100").Value = 8

# --- Rows 123/124: Hong Kong / Tunez swapped labels + refreshed numbers ---
$ws.Range("A123").Value = "Tunez"
$ws.Range("B123").Value = 1125
$ws.Range("C123").Value = 15
$ws.Range("D123").Value = 1002
$ws.Range("E123").Value = 74
$ws.Range("H123").Value = 49

$ws.Range("A124").Value = "Hong Kong"
$ws.Range("B124").Value = 1113
$ws.Range("D124").Value = 1069
$ws.Range("E124").Value = 40
$ws.Range("H124").Value = 4

# --- Rows 206/207: Groenlandia / Islas Malvinas swapped labels ---
$ws.Range("A206").Value = "Islas Malvinas"
$ws.Range("A207").Value = "Groenlandia"

# --- Rows 208/209: Santa Sede / Islas Turcas y Caicos swapped labels + data ---
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0

# --- Rows 213/214: Papua Nueva Guinea / Islas Virgenes Britanicas swapped labels + data ---
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
